# Insert a new product row (row 73: "مجموعه برد") into the pharmacy
# inventory sheet, pushing the totals row and the footer row down by one,
# and bump the grand-total cell by the new row's quantity.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new blank row at 73 (shifts old 73->74, 74->75, and
#        all existing mergeCells below it automatically). ---------------
$ws.Rows.Item(73).Insert()

# --- 2. Clone the formatting of the row above (row 72, the last
#        existing product row) into the new row 73. Copy per column-group
#        so each group keeps its own style id (A | B:G | H:K | L:M | N),
#        matching the layout used by every other product row. ----------
$ws.Range("A72:K72").Copy()
$ws.Range("A73:K73").PasteSpecial(-4122)
$ws.Range("L72:M72").Copy()
$ws.Range("L73:M73").PasteSpecial(-4122)
$ws.Range("N72").Copy()
$ws.Range("N73").PasteSpecial(-4122)

# Re-create the merged cells for the new row (Insert only carries the
# merges that were already below the insertion point).
$ws.Range("B73:G73").Merge()
$ws.Range("H73:K73").Merge()
$ws.Range("L73:M73").Merge()

# Row heights: new product row + new totals row use the standard 25.5
# height, the footer row settles at 16.5 (was 17.25 before the insert).
$ws.Rows.Item(73).RowHeight = 25.5
$ws.Rows.Item(74).RowHeight = 25.5
$ws.Rows.Item(75).RowHeight = 16.5

# --- 3. Fill in the new product row's data. -----------------------------
$ws.Range("A73").Value = 70
$ws.Range("B73").Value = "مجموعه برد"
$ws.Range("H73").Value = "-1:0"
$ws.Range("L73").Value = 8
$ws.Range("N73").Value = 1

# --- 4. Bump the grand total (now on row 74) by the new row's quantity. -
$total = $ws.Range("K74").Value()
$ws.Range("K74").Value = $total + 8
